# Apply updated cryptos list values per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.589.32'
$ws.Range('E2').Value = '  +0.74%  '
$ws.Range('D3').Value = '3.390.19'
$ws.Range('E3').Value = '  +0.61%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.43'
$ws.Range('E5').Value = '  +0.68%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.91'
$ws.Range('E6').Value = '  +0.67%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.389.04'
$ws.Range('E8').Value = '  +0.60%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.474'
$ws.Range('E9').Value = '  -0.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.47'
$ws.Range('E10').Value = '  -1.76%  '
$ws.Range('E11').Value = '  +2.23%  '
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('D13').Value = '3.966.40'
$ws.Range('E13').Value = '  +0.50%  '
$ws.Range('E15').Value = '  +1.85%  '
$ws.Range('D16').Value = '3.384.75'
$ws.Range('E16').Value = '  +0.38%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.78'
$ws.Range('E17').Value = '  +2.34%  '
$ws.Range('D18').Value = '61.639.22'
$ws.Range('E18').Value = '  +0.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.20'
$ws.Range('E19').Value = '  +1.92%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.83'
$ws.Range('E20').Value = '  +0.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.49'
$ws.Range('E21').Value = '  +0.56%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '379.00'
$ws.Range('E22').Value = '  +1.30%  '
$ws.Range('E23').Value = '  -1.36%  '
$ws.Range('D24').Value = '3.525.53'
$ws.Range('E24').Value = '  +0.54%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '71.14'
$ws.Range('E26').Value = '  +0.97%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000125'
$ws.Range('E27').Value = '  +6.64%  '
$ws.Range('E28').Value = '  -0.12%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.55'
$ws.Range('E29').Value = '  -1.92%  '
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('E31').Value = '  +3.59%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.18'
$ws.Range('E32').Value = '  +0.66%  '
$ws.Range('E33').Value = '  +1.11%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.38'
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.34'
$ws.Range('E37').Value = '  +0.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.84'
$ws.Range('E38').Value = '  -1.25%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '164.82'
$ws.Range('E39').Value = '  +0.84%  '
$ws.Range('E40').Value = '  -0.63%  '
$ws.Range('E41').Value = '  +2.47%  '
$ws.Range('B42').Value = 'Mantle'
$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.781'
$ws.Range('E42').Value = '  +2.68%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.999'
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('E44').Value = '  +7.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '25.29'
$ws.Range('E45').Value = '  +9.60%  '
$ws.Range('E46').Value = '  -0.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '41.37'
$ws.Range('E47').Value = '  +0.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.86'
$ws.Range('E48').Value = '  -1.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '22.69'
$ws.Range('E49').Value = '  -1.42%  '
$ws.Range('D50').Value = '2.350.08'
$ws.Range('E50').Value = '  +6.21%  '
$ws.Range('E51').Value = '  +2.09%  '
